$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.70%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.81%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.085"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.23%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07800"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.10%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.264"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.73%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.083"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.51%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.047"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.97%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9291"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.21%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1833"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.44%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08999"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-11.69%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08527"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.91%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03774"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "13.28%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09936"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.44%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.21%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005694"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.47%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.17%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.185"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.28%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.96%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1322"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.55%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.587"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.06%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.51%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.44%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.51%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004530"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.34%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.02%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-19.86%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.64%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04730"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.60%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007947"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.06%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1418"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.97%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007993"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-18.09%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002306"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009625"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.44%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006223"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.14%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.98%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "91.61%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002696"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.82%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.98%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.98%"
